# Update "想去人数" (F column) figures on the 展览 sheet and the
# corresponding rows on the 全部类型 sheet, mirroring the upstream
# gh-pages data refresh.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new F value on "展览"
$exhibitUpdates = @{
    6  = 80
    7  = 563
    8  = 74
    9  = 6828
    10 = 160
    13 = 178
    15 = 1105
    16 = 16241
    17 = 1597
    18 = 42
    22 = 11389
    24 = 1034
    25 = 4484
    26 = 330
    28 = 49
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new F value on "全部类型"
$allUpdates = @{
    6  = 80
    7  = 563
    9  = 74
    10 = 6828
    11 = 160
    14 = 178
    17 = 1105
    18 = 16241
    19 = 1597
    20 = 42
    26 = 11389
    28 = 1034
    29 = 4484
    30 = 330
    32 = 49
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
